# Add a "withdrawal report" to the "Foreign Currencies" sheet replacing the
# previous forex sell-order rows with updated buy/sell-date breakdown rows,
# and update the dependent summary figure on "ELSTER - Summary".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foreign Currencies")

# --- Row 4: keep Buy Date (2022-09-05) / Sell Date (2022-10-12), update amounts
$ws.Cells.Item(4, 2).Value = 2582.03
$ws.Cells.Item(4, 7).Value = 57.39

# --- Row 5: becomes Buy Date 2022-09-05 / Sell Date 2022-12-01
$ws.Cells.Item(5, 2).Value = 849.9400000000001
$ws.Cells.Item(5, 3).Value = "'2022-09-05"
$ws.Cells.Item(5, 3).Style = "Normal"
$ws.Cells.Item(5, 4).Value = "'2022-12-01"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 6).Value = 0.96
$ws.Cells.Item(5, 7).Value = -43.77

# --- Row 6: becomes Buy Date 2022-09-22 / Sell Date 2022-12-01
$ws.Cells.Item(6, 2).Value = 150.06
$ws.Cells.Item(6, 3).Value = "'2022-09-22"
$ws.Cells.Item(6, 3).Style = "Normal"
$ws.Cells.Item(6, 4).Value = "'2022-12-01"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = 1.01
$ws.Cells.Item(6, 6).Value = 0.96
$ws.Cells.Item(6, 7).Value = -8.279999999999999

# --- Remove the now-obsolete detail rows (old rows 7-11); the summary rows
# that used to be rows 12-15 shift up to become rows 7-10.
$ws.Rows.Item(7).Resize(5).EntireRow.Delete()

# --- Update the summary totals (new rows 8-10) for the recomputed figures.
$ws.Cells.Item(8, 7).Value = 22.05
$ws.Cells.Item(9, 7).Value = 74.09999999999999
$ws.Cells.Item(10, 7).Value = -52.05

# --- Propagate the updated "Gains (incl. losses)" total to the ELSTER summary sheet.
$ws6 = $wb.Worksheets.Item("ELSTER - Summary")
$ws6.Cells.Item(7, 3).Value = 22.05
